$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.697.25'
$ws.Range('E2').Value = '  +0.77%  '
$ws.Range('D3').Value = '2.359.32'
$ws.Range('E3').Value = '  +4.84%  '
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.659'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.07%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '234.61'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.68%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '73.81'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +14.94%  '
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.528'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +20.42%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0982'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.02%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '27.32'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.81%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.107'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.03%  '
$ws.Range('D13').Value = '2.710.82'
$ws.Range('E13').Value = '  +4.86%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '16.54'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +10.90%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.66'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +10.60%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.879'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +7.11%  '
$ws.Range('D17').Value = '2.370.35'
$ws.Range('E17').Value = '  +5.37%  '
$ws.Range('D18').Value = '43.692.49'
$ws.Range('E18').Value = '  +0.94%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0000101'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.50%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.46'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +6.86%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '75.62'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.50%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '251.07'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.85%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.84'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.09%  '
$ws.Range('E24').Value = '  -0.01%  '
$ws.Range('E25').Value = '  +1.93%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.22'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +5.38%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.25'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.07%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '22.50'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.41%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '172.11'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.95%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.54'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +7.59%  '
$ws.Range('E31').Value = '  +1.87%  '
$ws.Range('E32').Value = '  +4.83%  '
$ws.Range('E33').Value = '  +3.66%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0702'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.71%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.11'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.69%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.75'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.01%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.60'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.75%  '
$ws.Range('E38').Value = '  +7.47%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0263'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +5.29%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '19.53'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +13.27%  '
$ws.Range('E41').Value = '  -0.10%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.90'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.69%  '
$ws.Range('E43').Value = '  +9.73%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '99.90'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.47%  '
$ws.Range('B45').Value = 'Cronos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0969'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.56%  '
$ws.Range('B46').Value = 'TrustWalletToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.21'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.74%  '
$ws.Range('B47').Value = 'FTXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.44'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.11%  '
$ws.Range('E48').Value = '  +11.02%  '
$ws.Range('D49').Value = '1.442.36'
$ws.Range('E49').Value = '  +1.14%  '
$ws.Range('B50').Value = 'TerraClassic'
$ws.Range('C50').Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.000205'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.97%  '
$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').Value = '2.586.12'
$ws.Range('E51').Value = '  +4.58%  '
